# Generate Report for Handoff
# The a.md file finished its handoff/handback cycle (now occupies row 2 on each
# per-locale sheet) while the 747f25f9...md file is kicked back to "Ready for
# handoff" with a stale-handback error (now occupies row 3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 (a.md) status/date refreshed; row 2 (747f25f9) unchanged.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-02-28 07:13:52"

# Hyperlink display text for B2/B3 is swapped (the underlying link targets stay put).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/4dd14f3d77088b53b957fe91238fa8f13fda16f9/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md", "", "", "e2e\a.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2e44ee994f989fe5ba6bbf2034a32df2b096a305/e2e/a.md", "", "", "e2e\747f25f9-2892-47dc-87c9-7c7b9ba732bf.md")

# ---------------------------------------------------------------------------
# zh-cn sheet: row 2 becomes the a.md entry, row 3 becomes the 747f25f9 entry.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "a.md"
$wsZhCn.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-02-28 07:06:51"
$wsZhCn.Range("I2").Value = "TestHandoff_201702280307"
$wsZhCn.Range("J2").Value = "a.md"
$wsZhCn.Range("K2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("L2").Value = "2017-02-28 07:08:33"
$wsZhCn.Range("M2").Value = "TestHandback_201702280308"

$wsZhCn.Range("A3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.4ffb29e4e7febafd8e434a57081f966f0cf01f60.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2017-02-28 07:13:36"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md"
$wsZhCn.Range("K3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.4ffb29e4e7febafd8e434a57081f966f0cf01f60.zh-cn.xlf"
$wsZhCn.Range("L3").Value = "2017-02-28 07:12:29"
$wsZhCn.Range("M3").Value = ""
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/4dd14f3d77088b53b957fe91238fa8f13fda16f9/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/c3d721df7fa3252efa36f801696e79fcd167a4f4/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md."

# Hyperlinks A2/J2/A3/J3: swap the display text, keep original link targets.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/4dd14f3d77088b53b957fe91238fa8f13fda16f9/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/ebb99c9518028e564908f869bb19180fc4f8be78/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2e44ee994f989fe5ba6bbf2034a32df2b096a305/e2e/a.md", "", "", "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/d73d9e7d25161b44a47d361a3be07ca5c65aad38/e2e/a.md", "", "", "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md")

# Column R (Error Detail) widened to fit the new long error message.
$wsZhCn.Columns.Item(18).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: row 2 becomes the a.md entry, row 3 becomes the 747f25f9 entry.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "a.md"
$wsDeDe.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-02-28 07:07:10"
$wsDeDe.Range("I2").Value = "TestHandoff_201702280307"
$wsDeDe.Range("J2").Value = "a.md"
$wsDeDe.Range("K2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("L2").Value = "2017-02-28 07:08:55"
$wsDeDe.Range("M2").Value = "TestHandback_201702280308"

$wsDeDe.Range("A3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.4ffb29e4e7febafd8e434a57081f966f0cf01f60.de-de.xlf"
$wsDeDe.Range("H3").Value = "2017-02-28 07:13:52"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md"
$wsDeDe.Range("K3").Value = "747f25f9-2892-47dc-87c9-7c7b9ba732bf.4ffb29e4e7febafd8e434a57081f966f0cf01f60.de-de.xlf"
$wsDeDe.Range("L3").Value = "2017-02-28 07:12:51"
$wsDeDe.Range("M3").Value = ""
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/4dd14f3d77088b53b957fe91238fa8f13fda16f9/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/c3d721df7fa3252efa36f801696e79fcd167a4f4/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md."

# Hyperlinks A2/J2/A3/J3: swap the display text, keep original link targets.
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/4dd14f3d77088b53b957fe91238fa8f13fda16f9/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b804ed9095001b52b18796375317b2a9240d27a4/e2e/747f25f9-2892-47dc-87c9-7c7b9ba732bf.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2e44ee994f989fe5ba6bbf2034a32df2b096a305/e2e/a.md", "", "", "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/c51b84c664873af828173949009d24952370dfb9/e2e/a.md", "", "", "747f25f9-2892-47dc-87c9-7c7b9ba732bf.md")

# Column R (Error Detail) widened to fit the new long error message.
$wsDeDe.Columns.Item(18).ColumnWidth = 39.17
